$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# Helper: replace the contents of a whole paragraph with a specific,
# hand-built sequence of <w:r> runs (keeping the paragraph's own
# <w:p ...> attributes, e.g. w14:paraId / rsid, untouched).
# ---------------------------------------------------------------------
function Set-ParagraphRuns($paragraph, [string]$runsXml) {
    $pRange = $paragraph.Range

    # Pull the paragraph's original opening-tag attributes (paraId,
    # textId, rsid...) out of its flat WordprocessingML so the
    # rebuilt <w:p> keeps them instead of reverting to a bare <w:p>.
    $owx = $pRange.WordOpenXML
    $pOpenTag = "<w:p>"
    if ($owx -match "<w:p(\s[^>]*)?>") {
        $pOpenTag = $matches[0]
    }

    $pkg = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
           '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
           '<pkg:xmlData>' +
           '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml">' +
           '<w:body>' + $pOpenTag + $runsXml + '</w:p></w:body>' +
           '</w:document>' +
           '</pkg:xmlData></pkg:part></pkg:package>'

    # The Range must span the *whole* paragraph (including its mark)
    # for InsertXML to replace the paragraph's contents in place
    # instead of appending a sibling paragraph.
    $pRange.InsertXML($pkg)
}

function Get-ParagraphByText($doc, [string]$text) {
    foreach ($p in $doc.Paragraphs) {
        # Paragraph.Range.Text includes the trailing paragraph-mark
        # character (CR, 0x0D) - strip it before comparing.
        $t = $p.Range.Text.TrimEnd([char]13)
        if ($t -eq $text) { return $p }
    }
    return $null
}

# --- "Tableau: 40 points" -> "Tableau: 45 points" ----------------------
# split as: "Tableau: 4" + "5" + " points"
$target1 = Get-ParagraphByText $d "Tableau: 40 points"
$runs1 = '<w:r><w:t>Tableau: 4</w:t></w:r>' +
         '<w:r><w:t>5</w:t></w:r>' +
         '<w:r><w:t xml:space="preserve"> points</w:t></w:r>'
Set-ParagraphRuns $target1 $runs1

# --- "SQL:30 points" -> "SQL:25 points" ---------------------------------
# split as: "SQL:" + "25" + " points"
$target2 = Get-ParagraphByText $d "SQL:30 points"
$runs2 = '<w:r><w:t>SQL:</w:t></w:r>' +
         '<w:r><w:t>25</w:t></w:r>' +
         '<w:r><w:t xml:space="preserve"> points</w:t></w:r>'
Set-ParagraphRuns $target2 $runs2
